# Auto update stock data
# The "Date_1" column (column A) contains the as-of date for each company's
# most recent risk-score snapshot, stored as plain text (e.g. "2026/01/18").
# This refresh bumps that snapshot date forward by one day, for every block
# of data in the sheet (rows 2, 8, 14, 20, ... i.e. every 6th row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$targetRows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)
$oldDate = "2026/01/18"
$newDate = "2026/01/19"

foreach ($r in $targetRows) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldDate) {
        # Force a Text number format so Excel keeps the value as the literal
        # string "2026/01/19" instead of auto-converting it into a date
        # serial number (the source cell holds plain text, not a date).
        $cell.NumberFormat = "@"
        $cell.Value = $newDate
    }
}
